$d = $word.ActiveDocument

# Locate the last paragraph in the document (the existing bullet item
# "Which connections do we allow in the unit graph? ...") and append a
# new list item right after it, reusing the same list formatting
# (ListParagraph style + numPr ilvl 0 / numId 1).
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = "Need to specify how the data is stored and retrieved. We should be able to define a standard data interface at some high level so that we can compare data between different generators (and between identical generators). Also visualizers will need to use the data interface somehow"
